$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.870.77"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.814.38"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4658"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3684"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07364"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8702"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.816.45"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.365"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07064"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.504"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008688"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "26.894.57"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.336"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").Value = "2.047.45"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.902"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.178"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.319"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08926"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7659"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.165"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.504"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.903"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.087"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01960"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05281"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.930"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.248"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1660"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.423"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4923"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.671"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06286"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
